$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $ref, $value) {
    $ws.Range($ref).Value = $value
}

function Clear-Cell($ws, $ref) {
    $ws.Range($ref).ClearContents()
}

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
Set-Cell $ws "H33" 1034.2858
Set-Cell $ws "I33" 174.70589
Set-Cell $ws "K33" 174.70589
Set-Cell $ws "M33" 54.29410999999999
# Row 132
Set-Cell $ws "H132" 2095.913
Set-Cell $ws "I132" 1940.4667
Set-Cell $ws "K132" 5821.4001
Set-Cell $ws "M132" -3291.4001
# Row 134
Set-Cell $ws "H134" 69000
Set-Cell $ws "J134" 69000
Set-Cell $ws "L134" 69000
Set-Cell $ws "N134" -79140
# Row 136
Set-Cell $ws "H136" 30000
Set-Cell $ws "J136" 0
Set-Cell $ws "L136" 0
Clear-Cell $ws "N136"
# Row 137
Set-Cell $ws "H137" 2599.48
Set-Cell $ws "I137" 1870
Set-Cell $ws "K137" 5610
Set-Cell $ws "M137" -3060
# Row 138
Set-Cell $ws "H138" 3523.5908
Set-Cell $ws "I138" 3009.75
Set-Cell $ws "J138" 3817.2144
Set-Cell $ws "K138" 9029.25
Set-Cell $ws "L138" 11451.6432
Set-Cell $ws "M138" -3889.25
Set-Cell $ws "N138" -21731.6432
# Row 139
Set-Cell $ws "H139" 69996.5
Set-Cell $ws "J139" 69996.5
Set-Cell $ws "L139" 69996.5
Set-Cell $ws "N139" -80276.5
# Row 140
Set-Cell $ws "H140" 70490
Set-Cell $ws "J140" 70490
Set-Cell $ws "L140" 70490
Set-Cell $ws "N140" -80850

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 15
Set-Cell $ws "H15" 10.666667
Set-Cell $ws "I15" 10.666667
Set-Cell $ws "K15" 10.666667
Set-Cell $ws "M15" 339.333333
# Row 32
Set-Cell $ws "H32" 2535.94
Set-Cell $ws "I32" 1990.8889
Set-Cell $ws "J32" 7441.4
Set-Cell $ws "K32" 1990.8889
Set-Cell $ws "L32" 7441.4
Set-Cell $ws "M32" -1703.8889
Set-Cell $ws "N32" -8015.4
# Row 45
Set-Cell $ws "H45" 76926010
Set-Cell $ws "I45" 100001520
Set-Cell $ws "J45" 7635.3335
Set-Cell $ws "K45" 100001520
Set-Cell $ws "L45" 7635.3335
Set-Cell $ws "M45" -100001143
Set-Cell $ws "N45" -8389.333500000001
# Row 61
Set-Cell $ws "H61" 4420.9673
Set-Cell $ws "I61" 3434.7632
Set-Cell $ws "K61" 3434.7632
Set-Cell $ws "M61" -3222.7632
# Row 74
Set-Cell $ws "H74" 5295066
Set-Cell $ws "I74" 7096219
Set-Cell $ws "K74" 7096219
Set-Cell $ws "M74" -7095345
# Row 77
Set-Cell $ws "H77" 5295066
Set-Cell $ws "I77" 7096219
Set-Cell $ws "K77" 35481095
Set-Cell $ws "M77" -35476727
# Row 86
Set-Cell $ws "H86" 0
Set-Cell $ws "I86" 0
Set-Cell $ws "K86" 0
Clear-Cell $ws "M86"
# Row 89
Set-Cell $ws "H89" 0
Set-Cell $ws "I89" 0
Set-Cell $ws "K89" 0
Clear-Cell $ws "M89"
# Row 132
Set-Cell $ws "H132" 4845.0557
Set-Cell $ws "I132" 2827.6428
Set-Cell $ws "K132" 8482.928400000001
Set-Cell $ws "M132" -5952.928400000001
# Row 134
Set-Cell $ws "H134" 73500
Set-Cell $ws "J134" 73500
Set-Cell $ws "L134" 73500
Set-Cell $ws "N134" -83640
# Row 135
Set-Cell $ws "H135" 65809.336
Set-Cell $ws "J135" 65809.336
Set-Cell $ws "L135" 65809.336
Set-Cell $ws "N135" -75949.336
# Row 136
Set-Cell $ws "H136" 4420.9673
Set-Cell $ws "I136" 3434.7632
Set-Cell $ws "K136" 10304.2896
Set-Cell $ws "M136" -7754.2896

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
Set-Cell $ws "H20" 1647.6
Set-Cell $ws "I20" 1526
Set-Cell $ws "J20" 1931.3334
Set-Cell $ws "K20" 1526
Set-Cell $ws "L20" 1931.3334
Set-Cell $ws "M20" -1279
Set-Cell $ws "N20" -2425.3334
# Row 107
Set-Cell $ws "H107" 2036.2858
Set-Cell $ws "I107" 1542.3334
Set-Cell $ws "K107" 1542.3334
Set-Cell $ws "M107" 377.6666
# Row 134
Set-Cell $ws "H134" 1974.5555
Set-Cell $ws "I134" 1221.6552
Set-Cell $ws "K134" 3664.9656
Set-Cell $ws "M134" -1129.9656
# Row 137
Set-Cell $ws "H137" 0
Set-Cell $ws "J137" 0
Set-Cell $ws "L137" 0
Clear-Cell $ws "N137"
# Row 138
Set-Cell $ws "H138" 58092.562
Set-Cell $ws "J138" 58092.562
Set-Cell $ws "L138" 58092.562
Set-Cell $ws "N138" -68372.56200000001
# Row 140
Set-Cell $ws "H140" 69880
Set-Cell $ws "J140" 69880
Set-Cell $ws "L140" 69880
Set-Cell $ws "N140" -80240

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
Set-Cell $ws "H31" 30739.244
Set-Cell $ws "I31" 2918.9
Set-Cell $ws "J31" 57234.81
Set-Cell $ws "K31" 2918.9
Set-Cell $ws "L31" 57234.81
Set-Cell $ws "M31" -2623.9
Set-Cell $ws "N31" -57824.81
# Row 34
Set-Cell $ws "H34" 30739.244
Set-Cell $ws "I34" 2918.9
Set-Cell $ws "J34" 57234.81
Set-Cell $ws "K34" 2918.9
Set-Cell $ws "L34" 57234.81
Set-Cell $ws "M34" -2716.9
Set-Cell $ws "N34" -57638.81
# Row 105
Set-Cell $ws "H105" 2958.4167
Set-Cell $ws "I105" 1914.4286
Set-Cell $ws "K105" 1914.4286
Set-Cell $ws "M105" -167.4286
# Row 122
Set-Cell $ws "H122" 5419.7827
Set-Cell $ws "I122" 2583.2307
Set-Cell $ws "J122" 9107.299999999999
Set-Cell $ws "K122" 7749.6921
Set-Cell $ws "L122" 27321.9
Set-Cell $ws "M122" -5299.6921
Set-Cell $ws "N122" -32221.9
# Row 134
Set-Cell $ws "H134" 3004.5
Set-Cell $ws "I134" 2513.5151
Set-Cell $ws "J134" 3776.0476
Set-Cell $ws "K134" 7540.5453
Set-Cell $ws "L134" 11328.1428
Set-Cell $ws "M134" -5005.5453
Set-Cell $ws "N134" -16398.1428
# Row 138
Set-Cell $ws "H138" 60780
Set-Cell $ws "J138" 60780
Set-Cell $ws "L138" 60780
Set-Cell $ws "N138" -71060

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 119
Set-Cell $ws "H119" 44687
Set-Cell $ws "I119" 57514.5
Set-Cell $ws "J119" 19032
Set-Cell $ws "K119" 172543.5
Set-Cell $ws "L119" 57096
Set-Cell $ws "M119" -167705.5
Set-Cell $ws "N119" -66772
# Row 131
Set-Cell $ws "H131" 8103644.5
Set-Cell $ws "I131" 2262.375
Set-Cell $ws "J131" 12154336
Set-Cell $ws "K131" 6787.125
Set-Cell $ws "L131" 36463008
Set-Cell $ws "M131" -1747.125
Set-Cell $ws "N131" -36473088

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
Set-Cell $ws "H122" 4846.1875
Set-Cell $ws "I122" 3802.818
Set-Cell $ws "K122" 11408.454
Set-Cell $ws "M122" -8958.454000000002
# Row 132
Set-Cell $ws "H132" 23176
Set-Cell $ws "I132" 32502.242
Set-Cell $ws "J132" 5072.1177
Set-Cell $ws "K132" 97506.726
Set-Cell $ws "L132" 15216.3531
Set-Cell $ws "M132" -94976.726
Set-Cell $ws "N132" -20276.3531
# Row 140
Set-Cell $ws "H140" 0
Set-Cell $ws "J140" 0
Set-Cell $ws "L140" 0
Clear-Cell $ws "N140"

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 35
Set-Cell $ws "H35" 1833.3334
Set-Cell $ws "I35" 1833.3334
Set-Cell $ws "J35" 0
Set-Cell $ws "K35" 1833.3334
Set-Cell $ws "L35" 0
Set-Cell $ws "M35" -1497.3334
Clear-Cell $ws "N35"
# Row 136
Set-Cell $ws "H136" 7449
Set-Cell $ws "I136" 3275
Set-Cell $ws "K136" 9825
Set-Cell $ws "M136" -7275
# Row 137
Set-Cell $ws "H137" 59992.668
Set-Cell $ws "J137" 69989
Set-Cell $ws "L137" 69989
Set-Cell $ws "N137" -80189
# Row 139
Set-Cell $ws "H139" 0
Set-Cell $ws "J139" 0
Set-Cell $ws "L139" 0
Clear-Cell $ws "N139"

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 14
Set-Cell $ws "H14" 0
Set-Cell $ws "J14" 0
Set-Cell $ws "L14" 0
Clear-Cell $ws "N14"
# Row 70
Set-Cell $ws "H70" 28499.75
Set-Cell $ws "J70" 28499.75
Set-Cell $ws "L70" 28499.75
Set-Cell $ws "N70" -29129.75
# Row 73
Set-Cell $ws "H73" 28499.75
Set-Cell $ws "J73" 28499.75
Set-Cell $ws "L73" 28499.75
Set-Cell $ws "N73" -30683.75
# Row 103
Set-Cell $ws "H103" 38333.332
Set-Cell $ws "J103" 38333.332
Set-Cell $ws "L103" 38333.332
Set-Cell $ws "N103" -40677.332
# Row 122
Set-Cell $ws "H122" 3999.9048
Set-Cell $ws "I122" 2944.111
Set-Cell $ws "K122" 8832.332999999999
Set-Cell $ws "M122" -6382.332999999999
# Row 132
Set-Cell $ws "H132" 4094.205
Set-Cell $ws "I132" 3618.1177
Set-Cell $ws "J132" 7331.6
Set-Cell $ws "K132" 10854.3531
Set-Cell $ws "L132" 21994.8
Set-Cell $ws "M132" -8324.3531
Set-Cell $ws "N132" -27054.8
# Row 135
Set-Cell $ws "H135" 61621.5
Set-Cell $ws "J135" 61621.5
Set-Cell $ws "L135" 61621.5
Set-Cell $ws "N135" -71761.5
# Row 136
Set-Cell $ws "H136" 2797.5557
Set-Cell $ws "I136" 2035.4762
Set-Cell $ws "J136" 13466.667
Set-Cell $ws "K136" 6106.4286
Set-Cell $ws "L136" 40400.001
Set-Cell $ws "M136" -3556.4286
Set-Cell $ws "N136" -45500.001
# Row 137
Set-Cell $ws "H137" 66495
Set-Cell $ws "J137" 66495
Set-Cell $ws "L137" 66495
Set-Cell $ws "N137" -76695
# Row 141
Set-Cell $ws "H141" 135674.17
Set-Cell $ws "J141" 135674.17
Set-Cell $ws "L141" 135674.17
Set-Cell $ws "N141" -146034.17
